$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Narrow columns C and D (15 -> 6, 10 -> 6) --------------------------
# ColumnWidth is expressed in "characters"; on save the engine adds a
# fixed 5/6 padding on top of that number to produce the OOXML `width`
# attribute, so back that padding out here to land exactly on width=6.
$ws.Columns.Item(3).ColumnWidth = 6 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 6 - (5/6)

# --- New "GLASS" manual-input block (rows 13-16) ------------------------
# This pushes the previous GRAND TOTAL (formerly E13/E14) down to E19/E20
# and updates the total to reflect the newly-added glass line item.
$ws.Range("E13").Value = "GLASS"
$ws.Range("F13").Value = "Total Glass1"

$ws.Range("E14").Value = "Part Number"
$ws.Range("F14").Value = "N/A"

$ws.Range("E15").Value = "Quantity"
$ws.Range("F15").Value = "83.75 units"

$ws.Range("E16").Value = "Price"
# Force text storage for the currency-looking string so it isn't
# auto-parsed into a number by Excel's input parser, then drop the
# temporary "Text" number format so no stray style lingers on the cell.
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "$83.75"
$ws.Range("F16").ClearFormats()

# --- Keep rows 17-18 blank but present (matches the source row layout) --
$ws.Rows.Item(17).OutlineLevel = 0
$ws.Rows.Item(18).OutlineLevel = 0

# --- Re-inserted GRAND TOTAL block (rows 19-20) --------------------------
$ws.Range("E19").Value = "GRAND TOTAL"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "$2238.23"
$ws.Range("E20").ClearFormats()
